$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest batch of tickers (data refresh dated 2024-04-02).
$ws.Cells.Item(449, 1).Value = "IMX-USD"
$ws.Cells.Item(450, 1).Value = "MNT-USD"
$ws.Cells.Item(451, 1).Value = "PEPE-USD"
$ws.Cells.Item(452, 1).Value = "GRT-USD"
$ws.Cells.Item(453, 1).Value = "TAO-USD"
